{"js": "// Remove the last paragraph of the document body\n// (\"Recall is the opposite. ... that they have cancer.\"),\n// which is the paragraph right before the final section break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.delete();\nawait context.sync();\n", "ps1": "# Remove the last paragraph of the document body\n# (\"Recall is the opposite. ... that they have cancer.\"),\n# which is the paragraph right before the final section break.\n\n$d = $word.ActiveDocument\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.Delete()\n"}
